# cv124142a.xlsx - "correção nos dados e inicio da analise PNAD 2009"
#
# The sheet had two stray "label-only" rows (no data) that were throwing
# off the alignment between the row labels in column A and the data in
# columns B:F ("situação do domicílio" at row 5 and "grandes regiões e
# unidades da federação" at row 8). Removing those two rows lets every
# remaining label (brasil / urbana / rural / norte / states / ...) line
# back up with its correct data row. The header row also had duplicate
# "unnamed: *_level_1" placeholder labels that get replaced with "total"
# to match the corrected header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two label-only rows (no B:F data). Delete the lower one
# first so the upper row's index doesn't shift before we get to it.
$ws.Rows(8).Delete()
$ws.Rows(5).Delete()

# Fix the second header row: the "unnamed: 1_level_1" / "unnamed:
# 5_level_1" placeholders become "total", matching the already-correct
# "total" label in C2.
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
